$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.525.35"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "3.058.26"
$ws.Range("E3").Value = "  +2.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "386.07"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.05"
$ws.Range("D6").ClearFormats()
$ws.Range("E7").Value = "  -0.58%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.79"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "3.548.24"
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.57"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "3.057.93"
$ws.Range("E16").Value = "  +1.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.971"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.65"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.07%  "
$ws.Range("D19").Value = "51.565.29"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.16"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.44"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.17"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.97"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.14"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.35%  "
$ws.Range("E26").Value = "  +4.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.83"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.28"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.170"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.73%  "
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.107"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.26"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.75"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.26%  "
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("E35").Value = "  -3.06%  "
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("E38").Value = "  +2.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.292"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +7.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.93"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.76%  "
$ws.Range("E41").Value = "  +1.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.56"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "125.24"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.46%  "
$ws.Range("E45").Value = "  +2.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.93"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.89%  "
$ws.Range("E47").Value = "  +3.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.43"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.02%  "
$ws.Range("D49").Value = "2.032.44"
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("D50").Value = "3.360.14"
$ws.Range("E50").Value = "  +2.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.206"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +6.64%  "
